# Add the "Record_Cooldown" sheet (cooldown module drop item module) after Property1.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Record_Cooldown"

# ---- Cell values ----
$ws2.Range("A1").Value = "Id"
$ws2.Range("B1").Value = "Cooldown"

$ws2.Range("A2").Value = "Row"
$ws2.Range("B2").Value = 8
$ws2.Range("A3").Value = "Col"
$ws2.Range("B3").Value = 2
$ws2.Range("A4").Value = "Public"
$ws2.Range("B4").Value = 0
$ws2.Range("A5").Value = "Private"
$ws2.Range("B5").Value = 1
$ws2.Range("A6").Value = "Save"
$ws2.Range("B6").Value = 0
$ws2.Range("A7").Value = "Cache"
$ws2.Range("B7").Value = 1
$ws2.Range("A8").Value = "Upload"
$ws2.Range("B8").Value = 0

$ws2.Range("A9").Value = "SkillID"
$ws2.Range("B9").Value = "Time"
$ws2.Range("A10").Value = "string"
$ws2.Range("B10").Value = "int"
$ws2.Range("A11").Value = "Desc"

# ---- Formatting: reuse existing look & feel from Property1 ----
# Body rows (A2:B8) look like Property1's interior columns (orange fill).
$ws1.Range("C2").Copy()
$ws2.Range("A2:B8").PasteSpecial(-4122)

# Header row (A1:B1): same font/border, red fill.
$ws1.Range("C2").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws1.Range("B2").Copy()
$ws2.Range("B1").PasteSpecial(-4122)
$ws2.Range("A1:B1").Interior.Color = 255

# Footer rows (A9:B11): same font/border, blue accent fill.
$ws1.Range("C2").Copy()
$ws2.Range("A9:B10").PasteSpecial(-4122)
$ws2.Range("A11").PasteSpecial(-4122)
$ws1.Range("B2").Copy()
$ws2.Range("B11").PasteSpecial(-4122)
$ws2.Range("A9:B11").Interior.ThemeColor = 5
$ws2.Range("B11").ClearContents()

$ws2.Columns.Item(2).ColumnWidth = 12.9

# ---- Data validations ----
$ws2.Range("A2:B3").Validation.Add(1, 1, 5, 0)

$ws2.Range("A10:B10").Validation.Add(3, 1, 1, '"int,string,float,object"')
$ws2.Range("C11:H11").Validation.Add(3, 1, 1, '"int,string,float,object"')
$ws2.Range("C1:H8").Validation.Add(3, 1, 1, '"int,string,float,object"')

$ws2.Range("A9:B9").Validation.Add(0)
$ws2.Range("A9:B9").Validation.IgnoreBlank = $false

$ws2.Range("A4:B6").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws2.Range("B7:B8").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# ---- View / selection: make the new sheet the active tab ----
$ws2.Range("D10").Select()

# ---- Defined name: residual AutoFilter database on the new sheet ----
$n = $ws2.Names.Add("_xlnm._FilterDatabase", "=Record_Cooldown!`$A`$1:`$B`$11")
$n.Visible = $false
